# Update NATMI TPM-derived metrics on the active worksheet (Ntn1-Dcc.xlsx)
# Reflects a refresh of the underlying TPM data: the "Ligand-expressing cells"
# count for the ECs/Ntn1 group moved from 2 to 3 (detection rate 2/3 -> 1),
# which in turn changes the ligand average/total expression, and all of the
# derived specificity and edge-weight metrics that depend on it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> FAPs (Ntn1 -> Dcc)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.644726333333334
$ws.Range("H2").Value = 4.934179
$ws.Range("I2").Value = 0.03084360558270512
$ws.Range("J2").Value = 0.03084360558270512
$ws.Range("M2").Value = 0.092904
$ws.Range("Q2").Value = 0.152801655272
$ws.Range("R2").Value = 1.375214897448
$ws.Range("S2").Value = 0.02952484364034644
$ws.Range("T2").Value = 0.02952484364034644

# Row 3: ECs -> MuSCs (Ntn1 -> Dcc)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.644726333333334
$ws.Range("H3").Value = 4.934179
$ws.Range("I3").Value = 0.03084360558270512
$ws.Range("J3").Value = 0.03084360558270512
$ws.Range("O3").Value = 0.04275641311851519
$ws.Range("P3").Value = 0.04275641311851518
$ws.Range("Q3").Value = 0.006825066041222224
$ws.Range("R3").Value = 0.061425594371
$ws.Range("S3").Value = 0.001318761942358682
$ws.Range("T3").Value = 0.001318761942358681

# Row 4: FAPs -> FAPs (Ntn1 -> Dcc)
$ws.Range("I4").Value = 0.828024694817689
$ws.Range("J4").Value = 0.828024694817689
$ws.Range("M4").Value = 0.092904
$ws.Range("S4").Value = 0.7926213288937315
$ws.Range("T4").Value = 0.7926213288937315

# Row 5: FAPs -> MuSCs (Ntn1 -> Dcc)
$ws.Range("I5").Value = 0.828024694817689
$ws.Range("J5").Value = 0.828024694817689
$ws.Range("O5").Value = 0.04275641311851519
$ws.Range("P5").Value = 0.04275641311851518
$ws.Range("S5").Value = 0.03540336592395758
$ws.Range("T5").Value = 0.03540336592395757

# Row 6: MuSCs -> FAPs (Ntn1 -> Dcc)
$ws.Range("I6").Value = 0.1411316995996059
$ws.Range("J6").Value = 0.1411316995996059
$ws.Range("M6").Value = 0.092904
$ws.Range("Q6").Value = 0.699177573528
$ws.Range("R6").Value = 6.292598161752001
$ws.Range("S6").Value = 0.135097414347407
$ws.Range("T6").Value = 0.135097414347407

# Row 7: MuSCs -> MuSCs (Ntn1 -> Dcc)
$ws.Range("I7").Value = 0.1411316995996059
$ws.Range("J7").Value = 0.1411316995996059
$ws.Range("O7").Value = 0.04275641311851519
$ws.Range("P7").Value = 0.04275641311851518
$ws.Range("S7").Value = 0.006034285252198935
$ws.Range("T7").Value = 0.006034285252198933
